$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 17: was blank, now filled in (yellow D, dark-blue/white E = N/A) ---
$ws.Range("E17").Value = "N/A"
$ws.Range("E17").Interior.Color = 6299648
$ws.Range("E17").Font.ThemeColor = 2

$ws.Range("D17").Value = "20 revisions: 17 insertions, 3 deletions"
$ws.Range("D17").Interior.Color = 65535

# --- Row 18: was blank, now filled in (yellow D and E) ---
$ws.Range("D18").Value = "78 revisions: 43 insertions, 35 deletions"
$ws.Range("D18").Interior.Color = 65535

$ws.Range("E18").Value = "3 revisions: 1 insertion, 2 deletions"
$ws.Range("E18").Interior.Color = 65535

# --- Text updates (values only; styles for these cells are unchanged) ---
$ws.Range("E10").Value = "3 revisions: 3 insertions, 0 deletions"

$ws.Range("D11").Value = "9 revisions: 7 insertions, 2 deletions"
$ws.Range("E11").Value = "1 revision: 1 insertion. 0 deletions"

$ws.Range("D12").Value = "43 revisions, 33 insertions, 10 deletions"
$ws.Range("E12").Value = "1 revision: 1 insertion. 0 deletions"

$ws.Range("D13").Value = "1 revision: 1 insertion. 0 deletions"

$ws.Range("D14").Value = "28 revisions: 23 insertions, 5 deletions"

$ws.Range("D15").Value = "39 revisions: 34 insertions, 5 deletions"
$ws.Range("E15").Value = "8 revisions: 8 insertions, 0 deletions"

$ws.Range("D16").Value = "11 revisions: 9 insertions, 2 deletions"

# --- Row 16 E: was blank, now filled in with a green "reviewed" status ---
$ws.Range("E16").Value = "review complete - no change needed"
$ws.Range("E16").Interior.Color = 5287936

# --- Selection moves from D16 to D17 ---
$ws.Range("D17").Select()
